$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed cell values (login page test data refactor)
$ws.Range("A2").Value = "james8928748234"
$ws.Range("B3").Value = "Tommy82379834893"
$ws.Range("A3").Value = "thomas798597241"

# Widen columns A and B to fit the new, longer values
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 19.833333333333332

# Move the active selection to A3
$ws.Range("A3").Select()
